$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FactorID")

# ---------------------------------------------------------------------------
# Append the 17 new rows (80-96) of factor-id metadata to the FactorID sheet.
# Columns: A = Chinese field name, B = factor id, C = english field name,
#          D = table id (column C carries the "wrap text" style already used
#          by the other rows in this sheet, e.g. C2).
# ---------------------------------------------------------------------------

$data = @(
    @(80, "截止日",              132001, "date",         132),
    @(81, "股数",                132002, "shares_hold",  132),
    @(82, "占总股本比例(%)",      132003, "ratio_hold",   132),
    @(83, "截止日",              131001, "date",          131),
    @(84, "股票代码",            131002, "stock_id",      131),
    @(85, "股票名称",            131003, "stock_name",    131),
    @(86, "买入金额",            131004, "amt_buy",       131),
    @(87, "卖出金额",            131005, "amt_sell",      131),
    @(88, "买入及卖出金额",      131006, "amt_trade",     131),
    @(89, "排名",                131007, "rank",          131),
    @(90, "截止日",              130001, "date",          130),
    @(91, "卖出成交额(元)",      130004, "amt_sell",      130),
    @(92, "买入成交额(元)",      130003, "amt_buy",       130),
    @(93, "买入及卖出成交额(元)", 130002, "amt_trade",    130),
    @(94, "买入成交数目",        130009, "vol_buy",       130),
    @(95, "卖出成交数目",        130010, "vol_sell",      130),
    @(96, "买入及卖出成交数目",  130011, "vol_trade",     130)
)

foreach ($row in $data) {
    $r = $row[0]
    $ws.Cells.Item($r, 1).Value = $row[1]
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

# Column C in the existing rows uses a wrap-text style (see C2); copy that
# formatting onto the newly added column-C cells so they match.
$ws.Range("C2").Copy()
$ws.Range("C80:C96").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Reflect the new extent / last touched cell as the active selection.
$ws.Range("F86").Select()
